$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You are working as a data analyst for a healthcare company. For a data analysis project, you need to concatenate the contents of cells A2, B2, and C2 into a single cell in Excel. You aim to create a cleaner presentation of the data.What formula should you use to concatenate these cells?",
        "ques_type": 2,
        "options": [
            "=TEXTJOIN(\",\", TRUE, A2, B2, C2)",
            "=CONCATENATE(A2 &amp B2 &amp C2)",
            "=A2 &amp B2 &amp C2",
            "=CONCATENATE(A2, B2, C2)"
        ],
        "score": "=TEXTJOIN(\",\", TRUE, A2, B2, C2)"
    },
    {
        "title": "You are a data analyst working with a complex Excel spreadsheet that involves various calculations.  However, when certain data is missing or erroneous, the formula generates errors that are difficult for other team members to interpret.What formula should you use to display a custom \"Data Error\" message in these instances?",
        "ques_type": 2,
        "options": [
            "=IFNA(formula, \"Data Error\")",
            "=ISERROR(formula, \"Data Error\")",
            "=IFERROR(formula, \"Data Error\")",
            "=ERROR(formula, \"Data Error\")"
        ],
        "score": "=IFERROR(formula, \"Data Error\")"
    },
    {
        "title": "You are a data analyst for an e-commerce company. You are working with a large dataset in Excel. However, the dataset contains inconsistent capitalization in the product names column, making it difficult to sort or filter the data. You must clean up the inconsistent capitalization. What function should you use to standardize the capitalization in this column?",
        "ques_type": 2,
        "options": [
            "=CLEAN(column_name) ",
            "=TRIM(column_name)",
            "=UPPER(column_name)",
            "=PROPER(column_name)"
        ],
        "score": "=PROPER(column_name)"
    },
    {
        "title": "You are a data analyst for a technology company. You are working with an extensive dataset in Excel related to product sales. You must use the VLOOKUP function to retrieve product prices from a separate table. You write the following formula:=VLOOKUP(A2, Table_Range, 2, FALSE)However, you encounter an error where the function is not returning the expected results. What action should you take to rectify this issue?",
        "ques_type": 2,
        "options": [
            "Change the VLOOKUP range_lookup parameter from FALSE to TRUE.",
            "Move the table_array (Table_Range) to the same worksheet.",
            "Decrease the VLOOKUP column index number from two to one.",
            "Check to see if the lookup value (A2) contains leading or trailing spaces."
        ],
        "score": "Check to see if the lookup value (A2) contains leading or trailing spaces."
    }
]
'@

# Stage the long multi-line text in a scratch cell far from the used range,
# then paste-special (values only) into A2 so the direct-assignment auto
# row-height side effect lands on the scratch row instead of row 2.
$ws.Range("Z100").Value2 = $newText
$ws.Range("Z100").Copy()
$ws.Range("A2").PasteSpecial(-4163)
$ws.Rows.Item(100).Delete()

# Remove row 1 (the old numeric/bordered cell). This shifts row 2 (our
# updated shared-string cell) up into row 1, with no special style.
$ws.Rows.Item(1).Delete()
